# Applies the "Add data for 2022-07-12" change:
#  - Rolls the "through July 03" header forward to "through July 04"
#    (sheet name + the corresponding shared-string/header cell)
#  - Adds a handful of new carjacking counts (mostly filling in previously
#    empty/zero cells with 1, plus two cells whose totals increase)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the sheet and update the "through July 0X" header text ---
$ws.Name = "Through 2022-07-04"
$ws.Range("B1").Value = "July 2022 (through July 04)"

# --- Row 2: Austin ---
$ws.Range("B2").Value = 1
$ws.Range("W2").Value = 1
$ws.Range("AK2").Value = 1
$ws.Range("AR2").Value = 2
$ws.Range("AY2").Value = 1

# --- Row 4: Auburn Gresham ---
$ws.Range("B4").Value = 1

# --- Row 5: Garfield Park ---
$ws.Range("AY5").Value = 1

# --- Row 6: Grand Crossing ---
$ws.Range("B6").Value = 3

# --- Row 8: North Lawndale ---
$ws.Range("AK8").Value = 1

# --- Row 11: Loop ---
$ws.Range("P11").Value = 1

# --- Row 14: West Pullman ---
$ws.Range("W14").Value = 1

# --- Row 44: New City ---
$ws.Range("AY44").Value = 1

# --- Row 49: Grand Boulevard ---
$ws.Range("I49").Value = 1

# --- Row 52: Chatham ---
$ws.Range("I52").Value = 1

# --- Row 58: Albany Park ---
$ws.Range("I58").Value = 1

# --- Row 62: Avondale ---
$ws.Range("G62").Value = 1
$ws.Range("I62").Value = 1
